$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.853.82"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "1.740.33"
$ws.Range("E3").Value = "  -0.85%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'223.80"
$ws.Range("E5").Value = "  -5.75%  "
$ws.Range("D7").Value = "'0.5136"
$ws.Range("E7").Value = "  +1.33%  "
$ws.Range("D8").Value = "'0.2787"
$ws.Range("E8").Value = "  +5.55%  "
$ws.Range("D9").Value = "'38.98"
$ws.Range("E9").Value = "  -4.14%  "
$ws.Range("D10").Value = "'0.06078"
$ws.Range("E10").Value = "  -2.23%  "
$ws.Range("D11").Value = "1.748.64"
$ws.Range("E11").Value = "  -0.33%  "
$ws.Range("D12").Value = "'0.06945"
$ws.Range("D13").Value = "'15.14"
$ws.Range("E13").Value = "  -2.61%  "
$ws.Range("D14").Value = "'0.6304"
$ws.Range("E14").Value = "  +3.88%  "
$ws.Range("D15").Value = "'4.476"
$ws.Range("E15").Value = "  +0.63%  "
$ws.Range("D16").Value = "'76.07"
$ws.Range("E16").Value = "  -2.32%  "
$ws.Range("D17").Value = "'1.000"
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("D19").Value = "25.861.68"
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("D20").Value = "'11.37"
$ws.Range("E20").Value = "  -2.56%  "
$ws.Range("D21").Value = "'0.000006559"
$ws.Range("E21").Value = "  -3.92%  "
$ws.Range("D22").Value = "1.962.27"
$ws.Range("E22").Value = "  -0.58%  "
$ws.Range("D23").Value = "'4.064"
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").Value = "'8.357"
$ws.Range("E24").Value = "  +2.20%  "
$ws.Range("D25").Value = "'5.088"
$ws.Range("E25").Value = "  -1.70%  "
$ws.Range("D26").Value = "'138.23"
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("D27").Value = "'1.500"
$ws.Range("E27").Value = "  +1.75%  "
$ws.Range("D28").Value = "'1.810"
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("E29").Value = "  -0.98%  "
$ws.Range("D30").Value = "'102.10"
$ws.Range("E30").Value = "  -0.81%  "
$ws.Range("D31").Value = "'0.08250"
$ws.Range("E31").Value = "  -0.14%  "
$ws.Range("D32").Value = "'3.600"
$ws.Range("E32").Value = "  -2.62%  "
$ws.Range("D33").Value = "'3.380"
$ws.Range("E33").Value = "  -0.62%  "
$ws.Range("D34").Value = "'0.04359"
$ws.Range("E34").Value = "  -0.49%  "
$ws.Range("D35").Value = "'2.624"
$ws.Range("E35").Value = "  -1.17%  "
$ws.Range("D36").Value = "'0.9626"
$ws.Range("E36").Value = "  -4.30%  "
$ws.Range("D37").Value = "'0.5997"
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("D38").Value = "'2.658"
$ws.Range("E38").Value = "  -1.66%  "
$ws.Range("D39").Value = "'0.01544"
$ws.Range("E39").Value = "  -0.39%  "
$ws.Range("D40").Value = "'1.000"
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("D41").Value = "'1.889"
$ws.Range("E41").Value = "  -3.22%  "
$ws.Range("D42").Value = "'100.08"
$ws.Range("E42").Value = "  -3.39%  "
$ws.Range("D43").Value = "'0.3800"
$ws.Range("E43").Value = "  -0.55%  "
$ws.Range("D44").Value = "'0.7180"
$ws.Range("E44").Value = "  -3.85%  "
$ws.Range("D45").Value = "'4.891"
$ws.Range("E45").Value = "  +0.68%  "
$ws.Range("D46").Value = "'0.05445"
$ws.Range("E46").Value = "  -0.63%  "
$ws.Range("D47").Value = "'6.238"
$ws.Range("E47").Value = "  +4.57%  "
$ws.Range("D48").Value = "'0.1093"
$ws.Range("E48").Value = "  +1.23%  "
$ws.Range("D49").Value = "'51.97"
$ws.Range("E49").Value = "  -0.14%  "
$ws.Range("D50").Value = "'29.51"
$ws.Range("E50").Value = "  -2.15%  "
$ws.Range("D51").Value = "'1.003"
$ws.Range("E51").Value = "  +0.00%  "
